# Auto-generated edit script: applies cell-value updates per the commit diff
# for Sheets/Ixion_Profits.xlsx (8 worksheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 368.83334
$ws.Cells.Item(2, 9).Value = 311.45456
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 311.45456
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = -198.45456
$ws.Cells.Item(2, 14).Value = -1226
$ws.Cells.Item(62, 8).Value = 13890432
$ws.Cells.Item(62, 9).Value = 18520000
$ws.Cells.Item(62, 10).Value = 1728.3334
$ws.Cells.Item(62, 11).Value = 18520000
$ws.Cells.Item(62, 12).Value = 1728.3334
$ws.Cells.Item(62, 13).Value = -18519376
$ws.Cells.Item(62, 14).Value = -2976.3334
$ws.Cells.Item(65, 8).Value = 13890432
$ws.Cells.Item(65, 9).Value = 18520000
$ws.Cells.Item(65, 10).Value = 1728.3334
$ws.Cells.Item(65, 11).Value = 92600000
$ws.Cells.Item(65, 12).Value = 8641.666999999999
$ws.Cells.Item(65, 13).Value = -92596880
$ws.Cells.Item(65, 14).Value = -14881.667
$ws.Cells.Item(116, 8).Value = 11208.75
$ws.Cells.Item(116, 9).Value = 21541
$ws.Cells.Item(116, 11).Value = 21541
$ws.Cells.Item(116, 13).Value = -18099
$ws.Cells.Item(129, 8).Value = 1038.2162
$ws.Cells.Item(129, 9).Value = 726.8461
$ws.Cells.Item(129, 10).Value = 1104.5737
$ws.Cells.Item(129, 11).Value = 2180.5383
$ws.Cells.Item(129, 12).Value = 3313.7211
$ws.Cells.Item(129, 13).Value = 2819.4617
$ws.Cells.Item(129, 14).Value = -13313.7211
$ws.Cells.Item(132, 8).Value = 1441
$ws.Cells.Item(132, 9).Value = 1164.4642
$ws.Cells.Item(132, 11).Value = 3493.3926
$ws.Cells.Item(132, 13).Value = -963.3925999999997
$ws.Cells.Item(135, 8).Value = 1432.1384
$ws.Cells.Item(135, 9).Value = 1002.63464
$ws.Cells.Item(135, 11).Value = 9023.71176
$ws.Cells.Item(135, 13).Value = -6488.71176
$ws.Cells.Item(137, 8).Value = 1495.1892
$ws.Cells.Item(137, 9).Value = 1037.16
$ws.Cells.Item(137, 10).Value = 2449.4167
$ws.Cells.Item(137, 11).Value = 3111.48
$ws.Cells.Item(137, 12).Value = 7348.250100000001
$ws.Cells.Item(137, 13).Value = -561.4800000000005
$ws.Cells.Item(137, 14).Value = -12448.2501
$ws.Cells.Item(138, 8).Value = 2594.2026
$ws.Cells.Item(138, 9).Value = 1206.3
$ws.Cells.Item(138, 10).Value = 3540.5
$ws.Cells.Item(138, 11).Value = 3618.9
$ws.Cells.Item(138, 12).Value = 10621.5
$ws.Cells.Item(138, 13).Value = 1521.1
$ws.Cells.Item(138, 14).Value = -20901.5
$ws.Cells.Item(141, 8).Value = 1557.5217
$ws.Cells.Item(141, 9).Value = 1199.6129
$ws.Cells.Item(141, 10).Value = 2297.2
$ws.Cells.Item(141, 11).Value = 3598.8387
$ws.Cells.Item(141, 12).Value = 6891.599999999999
$ws.Cells.Item(141, 13).Value = 1581.1613
$ws.Cells.Item(141, 14).Value = -17251.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1763.52
$ws.Cells.Item(32, 9).Value = 1560.5
$ws.Cells.Item(32, 10).Value = 3252.3333
$ws.Cells.Item(32, 11).Value = 1560.5
$ws.Cells.Item(32, 12).Value = 3252.3333
$ws.Cells.Item(32, 13).Value = -1273.5
$ws.Cells.Item(32, 14).Value = -3826.3333
$ws.Cells.Item(45, 8).Value = 17312.834
$ws.Cells.Item(45, 9).Value = 20575.4
$ws.Cells.Item(45, 11).Value = 20575.4
$ws.Cells.Item(45, 13).Value = -20198.4
$ws.Cells.Item(61, 8).Value = 4094.342
$ws.Cells.Item(61, 9).Value = 5395.24
$ws.Cells.Item(61, 10).Value = 1592.6154
$ws.Cells.Item(61, 11).Value = 5395.24
$ws.Cells.Item(61, 12).Value = 1592.6154
$ws.Cells.Item(61, 13).Value = -5183.24
$ws.Cells.Item(61, 14).Value = -2016.6154
$ws.Cells.Item(74, 8).Value = 1104.1154
$ws.Cells.Item(74, 9).Value = 952.35297
$ws.Cells.Item(74, 10).Value = 1390.7778
$ws.Cells.Item(74, 11).Value = 952.35297
$ws.Cells.Item(74, 12).Value = 1390.7778
$ws.Cells.Item(74, 13).Value = -78.35297000000003
$ws.Cells.Item(74, 14).Value = -3138.7778
$ws.Cells.Item(77, 8).Value = 1104.1154
$ws.Cells.Item(77, 9).Value = 952.35297
$ws.Cells.Item(77, 10).Value = 1390.7778
$ws.Cells.Item(77, 11).Value = 4761.76485
$ws.Cells.Item(77, 12).Value = 6953.889
$ws.Cells.Item(77, 13).Value = -393.7648500000005
$ws.Cells.Item(77, 14).Value = -15689.889
$ws.Cells.Item(132, 8).Value = 3086.6978
$ws.Cells.Item(132, 9).Value = 2060.138
$ws.Cells.Item(132, 10).Value = 5213.143
$ws.Cells.Item(132, 11).Value = 6180.414
$ws.Cells.Item(132, 12).Value = 15639.429
$ws.Cells.Item(132, 13).Value = -3650.414
$ws.Cells.Item(132, 14).Value = -20699.429
$ws.Cells.Item(136, 8).Value = 4094.342
$ws.Cells.Item(136, 9).Value = 5395.24
$ws.Cells.Item(136, 10).Value = 1592.6154
$ws.Cells.Item(136, 11).Value = 16185.72
$ws.Cells.Item(136, 12).Value = 4777.8462
$ws.Cells.Item(136, 13).Value = -13635.72
$ws.Cells.Item(136, 14).Value = -9877.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1327.1904
$ws.Cells.Item(94, 9).Value = 491.9091
$ws.Cells.Item(94, 10).Value = 2246
$ws.Cells.Item(94, 11).Value = 491.9091
$ws.Cells.Item(94, 12).Value = 2246
$ws.Cells.Item(94, 13).Value = -40.90910000000002
$ws.Cells.Item(94, 14).Value = -3148

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1346
$ws.Cells.Item(16, 9).Value = 1240
$ws.Cells.Item(16, 10).Value = 1452
$ws.Cells.Item(16, 11).Value = 1240
$ws.Cells.Item(16, 12).Value = 1452
$ws.Cells.Item(16, 13).Value = -953
$ws.Cells.Item(16, 14).Value = -2026
$ws.Cells.Item(31, 8).Value = 274591.88
$ws.Cells.Item(31, 9).Value = 1699.35
$ws.Cells.Item(31, 10).Value = 916691.9
$ws.Cells.Item(31, 11).Value = 1699.35
$ws.Cells.Item(31, 12).Value = 916691.9
$ws.Cells.Item(31, 13).Value = -1404.35
$ws.Cells.Item(31, 14).Value = -917281.9
$ws.Cells.Item(34, 8).Value = 274591.88
$ws.Cells.Item(34, 9).Value = 1699.35
$ws.Cells.Item(34, 10).Value = 916691.9
$ws.Cells.Item(34, 11).Value = 1699.35
$ws.Cells.Item(34, 12).Value = 916691.9
$ws.Cells.Item(34, 13).Value = -1497.35
$ws.Cells.Item(34, 14).Value = -917095.9
$ws.Cells.Item(58, 8).Value = 1451.4117
$ws.Cells.Item(58, 9).Value = 922.5925999999999
$ws.Cells.Item(58, 10).Value = 2046.3334
$ws.Cells.Item(58, 11).Value = 922.5925999999999
$ws.Cells.Item(58, 12).Value = 2046.3334
$ws.Cells.Item(58, 13).Value = -719.5925999999999
$ws.Cells.Item(58, 14).Value = -2452.3334
$ws.Cells.Item(99, 8).Value = 10429470
$ws.Cells.Item(99, 9).Value = 13842
$ws.Cells.Item(99, 10).Value = 31260724
$ws.Cells.Item(99, 11).Value = 13842
$ws.Cells.Item(99, 12).Value = 31260724
$ws.Cells.Item(99, 13).Value = -12344
$ws.Cells.Item(99, 14).Value = -31263720
$ws.Cells.Item(105, 8).Value = 1802.4138
$ws.Cells.Item(105, 9).Value = 1924.2858
$ws.Cells.Item(105, 10).Value = 1482.5
$ws.Cells.Item(105, 11).Value = 1924.2858
$ws.Cells.Item(105, 12).Value = 1482.5
$ws.Cells.Item(105, 13).Value = -177.2858000000001
$ws.Cells.Item(105, 14).Value = -4976.5
$ws.Cells.Item(113, 8).Value = 1346
$ws.Cells.Item(113, 9).Value = 1240
$ws.Cells.Item(113, 10).Value = 1452
$ws.Cells.Item(113, 11).Value = 1240
$ws.Cells.Item(113, 12).Value = 1452
$ws.Cells.Item(113, 13).Value = 930
$ws.Cells.Item(113, 14).Value = -5792
$ws.Cells.Item(126, 8).Value = 10429470
$ws.Cells.Item(126, 9).Value = 13842
$ws.Cells.Item(126, 10).Value = 31260724
$ws.Cells.Item(126, 11).Value = 41526
$ws.Cells.Item(126, 12).Value = 93782172
$ws.Cells.Item(126, 13).Value = -39056
$ws.Cells.Item(126, 14).Value = -93787112
$ws.Cells.Item(132, 8).Value = 2024.138
$ws.Cells.Item(132, 9).Value = 1606.3721
$ws.Cells.Item(132, 10).Value = 3221.7334
$ws.Cells.Item(132, 11).Value = 4819.1163
$ws.Cells.Item(132, 12).Value = 9665.200199999999
$ws.Cells.Item(132, 13).Value = -2289.1163
$ws.Cells.Item(132, 14).Value = -14725.2002
$ws.Cells.Item(134, 8).Value = 1773.7847
$ws.Cells.Item(134, 9).Value = 2141
$ws.Cells.Item(134, 10).Value = 1103.2174
$ws.Cells.Item(134, 11).Value = 6423
$ws.Cells.Item(134, 12).Value = 3309.6522
$ws.Cells.Item(134, 13).Value = -3888
$ws.Cells.Item(134, 14).Value = -8379.6522
$ws.Cells.Item(136, 8).Value = 1451.4117
$ws.Cells.Item(136, 9).Value = 922.5925999999999
$ws.Cells.Item(136, 10).Value = 2046.3334
$ws.Cells.Item(136, 11).Value = 2767.7778
$ws.Cells.Item(136, 12).Value = 6139.0002
$ws.Cells.Item(136, 13).Value = -217.7777999999998
$ws.Cells.Item(136, 14).Value = -11239.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1154378.6
$ws.Cells.Item(113, 10).Value = 455076.22
$ws.Cells.Item(113, 12).Value = 1365228.66
$ws.Cells.Item(113, 14).Value = -1369568.66
$ws.Cells.Item(137, 8).Value = 18600.586
$ws.Cells.Item(137, 9).Value = 8663.333000000001
$ws.Cells.Item(137, 10).Value = 29247.643
$ws.Cells.Item(137, 11).Value = 25989.999
$ws.Cells.Item(137, 12).Value = 87742.929
$ws.Cells.Item(137, 13).Value = -20889.999
$ws.Cells.Item(137, 14).Value = -97942.929

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1358.4333
$ws.Cells.Item(102, 9).Value = 1135.8572
$ws.Cells.Item(102, 10).Value = 1877.7778
$ws.Cells.Item(102, 11).Value = 1135.8572
$ws.Cells.Item(102, 12).Value = 1877.7778
$ws.Cells.Item(102, 13).Value = 486.1428000000001
$ws.Cells.Item(102, 14).Value = -5121.7778
$ws.Cells.Item(113, 8).Value = 50001230
$ws.Cells.Item(113, 9).Value = 142858200
$ws.Cells.Item(113, 11).Value = 142858200
$ws.Cells.Item(113, 13).Value = -142856030
$ws.Cells.Item(132, 8).Value = 2447.5898
$ws.Cells.Item(132, 9).Value = 2029.4
$ws.Cells.Item(132, 10).Value = 2887.7896
$ws.Cells.Item(132, 11).Value = 6088.200000000001
$ws.Cells.Item(132, 12).Value = 8663.3688
$ws.Cells.Item(132, 13).Value = -3558.200000000001
$ws.Cells.Item(132, 14).Value = -13723.3688

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 31252146
$ws.Cells.Item(40, 9).Value = 35715670
$ws.Cells.Item(40, 10).Value = 7475
$ws.Cells.Item(40, 11).Value = 35715670
$ws.Cells.Item(40, 12).Value = 7475
$ws.Cells.Item(40, 13).Value = -35715534
$ws.Cells.Item(40, 14).Value = -7747
$ws.Cells.Item(63, 8).Value = 22500
$ws.Cells.Item(63, 10).Value = 22500
$ws.Cells.Item(63, 12).Value = 22500
$ws.Cells.Item(63, 14).Value = -23998
$ws.Cells.Item(66, 8).Value = 22500
$ws.Cells.Item(66, 10).Value = 22500
$ws.Cells.Item(66, 12).Value = 67500
$ws.Cells.Item(66, 14).Value = -74988
$ws.Cells.Item(122, 8).Value = 2911059.8
$ws.Cells.Item(122, 9).Value = 3404327.5
$ws.Cells.Item(122, 10).Value = 1431257.1
$ws.Cells.Item(122, 11).Value = 10212982.5
$ws.Cells.Item(122, 12).Value = 4293771.300000001
$ws.Cells.Item(122, 13).Value = -10210532.5
$ws.Cells.Item(122, 14).Value = -4298671.300000001
$ws.Cells.Item(132, 8).Value = 10690686
$ws.Cells.Item(132, 9).Value = 12726155
$ws.Cells.Item(132, 10).Value = 4475.375
$ws.Cells.Item(132, 11).Value = 38178465
$ws.Cells.Item(132, 12).Value = 13426.125
$ws.Cells.Item(132, 13).Value = -38175935
$ws.Cells.Item(132, 14).Value = -18486.125
$ws.Cells.Item(137, 8).Value = 28825.8
$ws.Cells.Item(137, 10).Value = 28825.8
$ws.Cells.Item(137, 12).Value = 28825.8
$ws.Cells.Item(137, 14).Value = -39025.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 42300
$ws.Cells.Item(68, 10).Value = 42300
$ws.Cells.Item(68, 12).Value = 42300
$ws.Cells.Item(68, 14).Value = -43922
$ws.Cells.Item(71, 8).Value = 42300
$ws.Cells.Item(71, 10).Value = 42300
$ws.Cells.Item(71, 12).Value = 126900
$ws.Cells.Item(71, 14).Value = -135012
$ws.Cells.Item(122, 8).Value = 2348.4243
$ws.Cells.Item(122, 9).Value = 2312.4167
$ws.Cells.Item(122, 11).Value = 6937.250100000001
$ws.Cells.Item(122, 13).Value = -4487.250100000001
$ws.Cells.Item(132, 8).Value = 17017.273
$ws.Cells.Item(132, 9).Value = 19485.076
$ws.Cells.Item(132, 10).Value = 2484.6667
$ws.Cells.Item(132, 11).Value = 58455.228
$ws.Cells.Item(132, 12).Value = 7454.000100000001
$ws.Cells.Item(132, 13).Value = -55925.228
$ws.Cells.Item(132, 14).Value = -12514.0001
$ws.Cells.Item(136, 8).Value = 7814952
$ws.Cells.Item(136, 9).Value = 2580.561
$ws.Cells.Item(136, 11).Value = 7741.683000000001
$ws.Cells.Item(136, 13).Value = -5191.683000000001

